# Update the instructional banner text that appears under every section
# header ("2Decision Criteria", "3Decision Makers", ... "7Alternatives")
# throughout the pairwise_comp worksheet.
#
# The old text ("Enter judgments for the paiwise comparisons in the matrix
# or direct values in the green cells") is replaced everywhere by the new,
# reworded instructions. Excel stores this repeated text once in the
# shared-strings table, so updating every cell that used the old string to
# the new text causes the old (now-unreferenced) entry to drop out of the
# table and the new entry to be appended - matching the workbook's shared
# string table after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pairwise_comp")

$newInstructions = "Enter pairwise comparisons in the white cells of the table or numerical data in the green cells. For the Direct Values column, if the smallest value is best, invert the value before entering it (e.g., `$10 as =1/10) ."

$bannerCells = @(
    "A2", "A11", "A21", "A31", "A41", "A51", "A62", "A73", "A84", "A95",
    "A106", "A115", "A124", "A133", "A142", "A151", "A160", "A167", "A174",
    "A180", "A186", "A193", "A200", "A207", "A214", "A221"
)

foreach ($cellRef in $bannerCells) {
    $ws.Range($cellRef).Value = $newInstructions
}

# Update the view so the last edited table (row ~221) is the one in focus,
# matching the saved selection/scroll state from the author's last save.
$ws.Range("A221").Select()
